$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "D1_USD" (sheet1.xml): update B129, append rows 130:133
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("D1_USD")

# corrected value for the existing last row
$ws1.Range("B129").Value = 4.021039

# bring in formatting for the new rows by copying the row above, then
# overwrite with the real values/formulas
$ws1.Range("A129:E129").Copy()
$ws1.Range("A130:E130").PasteSpecial(-4122)
$ws1.Range("A131:E131").PasteSpecial(-4122)
$ws1.Range("A132:E132").PasteSpecial(-4122)

$ws1.Range("A130").Value = 45342
$ws1.Range("B130").Value = 4.0136430000000001
$ws1.Range("C130").Value = 4.027825

$ws1.Range("A131").Value = 45343
$ws1.Range("B131").Value = 3.9898500000000001
$ws1.Range("C131").Value = 4.0217590000000003

$ws1.Range("A132").Value = 45344
$ws1.Range("B132").Value = 3.9882620000000002
$ws1.Range("C132").Value = 3.9474418

$ws1.Range("D130:D132").Formula = "=B130-C130"
$ws1.Range("E130:E132").Formula = "=IF(D130<0,1,0)"

$ws1.Range("C133").Value = 4.0039740000000004

$ws1.Activate()
$ws1.Range("C134").Select()

# ---------------------------------------------------------------------------
# Sheet "D1_EUR" (sheet3.xml): append rows 404:407
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("D1_EUR")

$ws3.Range("A403:E403").Copy()
$ws3.Range("A404:E404").PasteSpecial(-4122)
$ws3.Range("A405:E405").PasteSpecial(-4122)
$ws3.Range("A406:E406").PasteSpecial(-4122)

$ws3.Range("A404").Value = 45342
$ws3.Range("B404").Value = 4.3247
$ws3.Range("C404").Value = 4.3487340000000003
$ws3.Range("D404").Formula = "=B404-C404"

$ws3.Range("A405").Value = 45343
$ws3.Range("B405").Value = 4.3125299999999998
$ws3.Range("C405").Value = 4.3379659999999998
$ws3.Range("D405").Formula = "=B405-C405"

$ws3.Range("A406").Value = 45344
$ws3.Range("B406").Value = 4.3152999999999997
$ws3.Range("C406").Value = 4.3258729999999996
$ws3.Range("D406").Formula = "=B406-C406"

$ws3.Range("E404:E406").Formula = "=IF(D404<0,1,0)"

$ws3.Range("C407").Value = 4.3224660000000004

$ws3.Activate()
$ws3.Range("D406").Select()

# ---------------------------------------------------------------------------
# Sheet "D5_EUR" (sheet5.xml): turn the "Nan" placeholders in B96:B99 into
# real numbers, then append rows 100:104
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("D5_EUR")

$ws5.Range("B96").Value = 4.3365299999999998
$ws5.Range("B97").Value = 4.3247
$ws5.Range("B98").Value = 4.3125299999999998
$ws5.Range("B99").Value = 4.3152999999999997

$ws5.Range("A99:C99").Copy()
$ws5.Range("A100:C100").PasteSpecial(-4122)
$ws5.Range("A101:C101").PasteSpecial(-4122)
$ws5.Range("A102:C102").PasteSpecial(-4122)
$ws5.Range("A103:C103").PasteSpecial(-4122)
$ws5.Range("A104:C104").PasteSpecial(-4122)

$ws5.Range("A100").Value = 45345
$ws5.Range("B100").Value = "Nan"
$ws5.Range("C100").Value = 4.3074820000000003

$ws5.Range("A101").Value = 45348
$ws5.Range("B101").Value = "Nan"
$ws5.Range("C101").Value = 4.3080907000000002

$ws5.Range("A102").Value = 45349
$ws5.Range("B102").Value = "Nan"
$ws5.Range("C102").Value = 4.3042490000000004

$ws5.Range("A103").Value = 45350
$ws5.Range("B103").Value = "Nan"
$ws5.Range("C103").Value = 4.3033146999999996

$ws5.Range("A104").Value = 45351
$ws5.Range("B104").Value = "Nan"
$ws5.Range("C104").Value = 4.3039784000000001

# ---------------------------------------------------------------------------
# Sheet "D1_OIL" (sheet6.xml): append rows 76:79
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("D1_OIL")

$ws6.Range("A75:E75").Copy()
$ws6.Range("A76:E76").PasteSpecial(-4122)
$ws6.Range("A77:E77").PasteSpecial(-4122)
$ws6.Range("A78:E78").PasteSpecial(-4122)

$ws6.Range("A76").Value = 45342
$ws6.Range("B76").Value = 78.269997000000004
$ws6.Range("C76").Value = 68.0334

$ws6.Range("A77").Value = 45343
$ws6.Range("B77").Value = 77.910004000000001
$ws6.Range("C77").Value = 67.810500000000005

$ws6.Range("A78").Value = 45344
$ws6.Range("B78").Value = 78.610000999999997
$ws6.Range("C78").Value = 68.6648

$ws6.Range("D76:D78").Formula = "=B76-C76"
$ws6.Range("E76:E78").Formula = "=D76/C76"

$ws6.Range("C79").Value = 69.136099999999999

$ws6.Activate()
$ws6.Range("C80").Select()

# ---------------------------------------------------------------------------
# D5_EUR ends up as the active tab/sheet, with B100 selected
# ---------------------------------------------------------------------------
$ws5.Activate()
$ws5.Range("B100").Select()
